$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab from "Sheet1" to "2016"
$ws.Name = "2016"

# 2. Fix the old/typo'd component name: ESP6288 -> ESP8266
$ws.Range("D7").Value = "Module wifi ESP8266"

# 3. Shrink the sheet-tab area a touch (tabRatio 991 -> 964)
$excel.Windows.Item(1).TabRatio = 0.964

# 4. A handful of previously-blank cells picked up the default "Normal" style
#    (same style already used by C1 / E1) after the user clicked/formatted
#    around that area of the sheet.
[void]$ws.Range("C1").Copy()
[void]$ws.Range("B1").PasteSpecial(-4122)
[void]$ws.Range("C1").Copy()
[void]$ws.Range("C6").PasteSpecial(-4122)
[void]$ws.Range("C1").Copy()
[void]$ws.Range("C7").PasteSpecial(-4122)
[void]$ws.Range("C1").Copy()
[void]$ws.Range("C8").PasteSpecial(-4122)
[void]$ws.Range("C1").Copy()
[void]$ws.Range("C9").PasteSpecial(-4122)
[void]$ws.Range("E1").Copy()
[void]$ws.Range("E9").PasteSpecial(-4122)

# 5. Move the saved cursor/selection position to D17
[void]$ws.Range("D17").Select()
